$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hours")

# Drop the stray formatted-but-empty last row carried over from the source file
$ws.Rows.Item(1048576).Delete()

# Insert a new column E (shifts the old "jobstr" column, and its data, to F)
$ws.Columns.Item(5).Insert()

# New header for the inserted column
$ws.Range("E1").Value = "full_pcnt"

# New percentage values for the inserted column, formatted to 3 decimals, centered
$pcnt = @(0.6, 0.625, 0.65, 0.6, 0.625, 0.65, 0.65, 0.65)
for ($i = 0; $i -lt $pcnt.Length; $i++) {
    $row = $i + 2
    $ws.Range("E$row").Value = $pcnt[$i]
}
$ws.Range("E2:E9").NumberFormat = "0.000"
$ws.Range("E2:E9").HorizontalAlignment = -4108

# Rename the job strings now living in column F
$jobstr = @("Capt G4", "Capt G3", "Capt G2", "F/O G4", "F/O G3", "F/O G2", "Capt G1", "F/O G1")
for ($i = 0; $i -lt $jobstr.Length; $i++) {
    $row = $i + 2
    $ws.Range("F$row").Value = $jobstr[$i]
}

# Column widths: A:F now share one uniform width, reverting to the default beyond F
$ws.Range("A1:F1").EntireColumn.ColumnWidth = 12.8542510121458

# Selection / view bookkeeping to match the saved state
$ws.Range("I15").Select()

# Cosmetic tab-split ratio (best effort; some hosts do not persist window chrome)
$wb.Windows.Item(1).TabRatio = 993
